$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that needs to move
# from 45186 (2023-09-17) to 45188 (2023-09-19) for every data row
# (rows 2 through 233).
for ($row = 2; $row -le 233; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
